# `excel`: better cell error handling
# Adds a new "cellerrors" worksheet (after "data types") that exercises
# several error-producing formulas (#DIV/0!, #NAME?) plus a defined name
# ("testname" -> cellerrors!$C$6) used by one of the formulas.

$wb = $excel.ActiveWorkbook

# --- add the new sheet as the LAST sheet (after "data types") ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "cellerrors"

# --- header row ----------------------------------------------------------
$ws.Range("A1").Value = "col1"
$ws.Range("B1").Value = "col 2"
$ws.Range("C1").Value = "column-3"

# --- row 2 -----------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Formula = "=100/-2"
$ws.Range("C2").Formula = "=SUM(A2:A6)"

# --- row 3 (errors: #DIV/0! and a #NAME? array formula) --------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Formula = "=100/0"
$ws.Range("C3").FormulaArray = "=4*te"

# --- row 4 (uses the defined name "testname", defined below) ---------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Formula = "=100/2"

# --- row 5 -------------------------------------------------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Formula = "=100/3"
$ws.Range("C5").Value = 3

# --- row 6 -------------------------------------------------------------------
$ws.Range("A6").Value = 5
$ws.Range("B6").Formula = "=100/4"
$ws.Range("C6").Value = 4

# --- workbook-scoped defined name: testname -> cellerrors!$C$6 -------------
$wb.Names.Add("testname", "=cellerrors!`$C`$6")

# now that the name exists, write the formula that references it
$ws.Range("C4").Formula = "=5*testname"

# --- page setup (matches exported worksheet's pageSetup element) -----------
$ws.PageSetup.Orientation = 1

# --- make "cellerrors" the active sheet / tab, with C4 selected ------------
$ws.Activate() | Out-Null
$ws.Range("C4").Select() | Out-Null
